# issue #5: stock data from json to db
#
# The json-to-db export now emits three extra columns that the other
# property sheets already had (category, source_file, index). Add them to
# the 股票 (stock) sheet:
#   - a "category" column right after "property_category"
#   - "source_file" and "index" columns appended at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at I for "category" (between property_category and
# date); the old I/J/K (date/legislator_name/legislator_id) shift right to
# J/K/L.
$ws.Columns("I:I").Insert()

# Header row
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data row
$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmp40191"
$ws.Range("N2").Value = 72

# Match the formatting of the rest of the header row (bold, centered,
# bordered) for the three newly-introduced header cells.
foreach ($addr in @("I1", "M1", "N1")) {
    $r = $ws.Range($addr)
    $r.Borders.LineStyle = 1
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
}
